# Update the "review" answers for rows 17-19 (G column) from "yes" to "no",
# and select the edited range (matching the author's active selection).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G17:G19").Value = "no"
$ws.Range("G17:G19").Select()
